# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ROCIO DEL CARMEN REDONDO MENDEZ (45432136) - periods renumbered sequentially
# starting at 2101, with the last period (2211) moved down to row 40 to make
# room for the two new WILLIAM JAVIER VANEGAS BARROS rows (33 and 35).
$ws.Range("E16").Value = "2101"
$ws.Range("F16").Value = 35112

$ws.Range("E17").Value = "2102"
$ws.Range("E18").Value = "2103"
$ws.Range("E19").Value = "2104"
$ws.Range("E20").Value = "2105"
$ws.Range("E21").Value = "2106"
$ws.Range("E22").Value = "2107"
$ws.Range("E23").Value = "2108"
$ws.Range("E24").Value = "2109"
$ws.Range("E25").Value = "2110"
$ws.Range("E26").Value = "2111"
# E27 stays "2112"

$ws.Range("E28").Value = "2201"
$ws.Range("E29").Value = "2202"
$ws.Range("E30").Value = "2203"
$ws.Range("E31").Value = "2204"
$ws.Range("E32").Value = "2205"

# New row 33: WILLIAM JAVIER VANEGAS BARROS, periodo 2205
$ws.Range("C33").Value = "1121333761"
$ws.Range("D33").Value = "WILLIAM JAVIER VANEGAS BARROS"
$ws.Range("E33").Value = "2205"
$ws.Range("F33").Value = 40000
$ws.Range("G33").Value = 1000000

# Row 34: ROCIO, periodo 2206
$ws.Range("E34").Value = "2206"

# New row 35: WILLIAM JAVIER VANEGAS BARROS, periodo 2206
$ws.Range("C35").Value = "1121333761"
$ws.Range("D35").Value = "WILLIAM JAVIER VANEGAS BARROS"
$ws.Range("E35").Value = "2206"
$ws.Range("F35").Value = 40000
$ws.Range("G35").Value = 1000000

$ws.Range("E36").Value = "2207"
$ws.Range("E37").Value = "2208"
$ws.Range("E38").Value = "2209"

# Row 39: now ROCIO, periodo 2210 (was ANGEL DURANGO MORALES 2211)
$ws.Range("C39").Value = "45432136"
$ws.Range("D39").Value = "ROCIO DEL CARMEN REDONDO MENDEZ"
$ws.Range("E39").Value = "2210"
$ws.Range("F39").Value = 35112
$ws.Range("G39").Value = 877803

# Row 40: now ROCIO, periodo 2211 (was WILLIAM JAVIER VANEGAS BARROS 2206)
$ws.Range("C40").Value = "45432136"
$ws.Range("D40").Value = "ROCIO DEL CARMEN REDONDO MENDEZ"
$ws.Range("E40").Value = "2211"
$ws.Range("F40").Value = 24578
$ws.Range("G40").Value = 877803

# Row 41: now ANGEL DURANGO MORALES, periodo 2211 (was WILLIAM JAVIER VANEGAS BARROS 2205)
$ws.Range("C41").Value = "73110179"
$ws.Range("D41").Value = "ANGEL DURANGO MORALES"
$ws.Range("E41").Value = "2211"
$ws.Range("F41").Value = 28000
